$wb = $excel.ActiveWorkbook

# --- Keywords sheet: add new keyword rows + widen column F, then reselect ---
$wsKeywords = $wb.Worksheets.Item("Keywords")

$keywordValues = @(
    "baseUrl",
    "startingPopUp",
    "inputMenu",
    "simpleMenu",
    "singleInput",
    "messageButton",
    "messageResult",
    "firstNumber",
    "sumButton",
    "result",
    "checkbox",
    "singleCheck",
    "singleResult",
    "singleResult",
    "multipleCheck1",
    "multipleCheck2",
    "multipleCheck3",
    "multipleCheck4",
    "multipleResult",
    "dropdown",
    "select",
    "monday",
    "tuesday",
    "wednesday",
    "thursday",
    "friday",
    "saturday",
    "sunday",
    "dayResult",
    "radio",
    "male",
    "female",
    "0-5",
    "5-15",
    "15",
    "radiobutton",
    "radioResult"
)

$row = 2
foreach ($val in $keywordValues) {
    $cell = $wsKeywords.Range("C$row")
    if ($val -eq "15") {
        # Force the plain numeric-looking keyword "15" to be stored as
        # text (it is a keyword name, not a number), matching the rest
        # of the column which is all shared-string text.
        $cell.NumberFormat = "@"
        $cell.Value = $val
        $cell.NumberFormat = "General"
    } else {
        $cell.Value = $val
    }
    $row = $row + 1
}

# New wide column for notes (column F), and an empty trailing row 39
# to match the extended used range.
$wsKeywords.Columns.Item(6).ColumnWidth = 65
$wsKeywords.Range("C39").NumberFormat = "General"

$wsKeywords.Range("B9").Select() | Out-Null

# --- Navigation sheet: move the remembered selection to B2 ---
$wsNavigation = $wb.Worksheets.Item("Navigation")
$wsNavigation.Range("B2").Select() | Out-Null

# --- SelectList sheet: move the remembered selection to C10 and keep it active ---
$wsSelectList = $wb.Worksheets.Item("SelectList")
$wsSelectList.Range("C10").Select() | Out-Null
$wsSelectList.Activate() | Out-Null
